$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.408.28'
$ws.Range("E2").Value = '  -1.55%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.186.48'
$ws.Range("E3").Value = '  -2.36%  '

$ws.Range("E4").Value = '  -0.42%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.64'
$ws.Range("E5").Value = '  +2.85%  '

$ws.Range("E6").Value = '  -0.95%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.81'
$ws.Range("E7").Value = '  -1.21%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("E9").Value = '  -5.23%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.43'
$ws.Range("E10").Value = '  -1.69%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0910'
$ws.Range("E11").Value = '  -1.49%  '

$ws.Range("E12").Value = '  +0.35%  '

$ws.Range("E13").Value = '  -2.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.510.77'
$ws.Range("E14").Value = '  -1.75%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.18'
$ws.Range("E15").Value = '  -3.65%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.168.73'
$ws.Range("E16").Value = '  -3.26%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.771'
$ws.Range("E17").Value = '  -5.36%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.349.74'
$ws.Range("E18").Value = '  -1.50%  '

$ws.Range("E19").Value = '  -3.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.70'
$ws.Range("E20").Value = '  -0.65%  '

$ws.Range("E21").Value = '  -2.29%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.70'
$ws.Range("E22").Value = '  -1.63%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.39'
$ws.Range("E23").Value = '  -10.81%  '

$ws.Range("E24").Value = '  -5.45%  '

$ws.Range("E25").Value = '  +0.01%  '

$ws.Range("E26").Value = '  -4.76%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.37'

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '38.47'
$ws.Range("E28").Value = '  +1.75%  '

$ws.Range("E29").Value = '  -0.38%  '

$ws.Range("E30").Value = '  -3.88%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '172.65'
$ws.Range("E31").Value = '  -0.71%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.04'
$ws.Range("E32").Value = '  -1.48%  '

$ws.Range("E33").Value = '  +3.30%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.13'
$ws.Range("E34").Value = '  -4.92%  '

$ws.Range("E35").Value = '  -1.26%  '

$ws.Range("E36").Value = '  -4.17%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0338'
$ws.Range("E37").Value = '  +1.92%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.23'
$ws.Range("E38").Value = '  -3.49%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.02'
$ws.Range("E39").Value = '  -9.26%  '

$ws.Range("E40").Value = '  -3.90%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.61'
$ws.Range("E41").Value = '  +12.29%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.16'
$ws.Range("E42").Value = '  -7.58%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '58.83'
$ws.Range("E43").Value = '  -2.40%  '

$ws.Range("E44").Value = '  -3.09%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.37'
$ws.Range("E45").Value = '  -3.83%  '

$ws.Range("E46").Value = '  -2.40%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.460'
$ws.Range("E47").Value = '  +3.34%  '

$ws.Range("E48").Value = '  -5.05%  '

$ws.Range("E49").Value = '  -1.67%  '

$ws.Range("E50").Value = '  -2.29%  '

$ws.Range("E51").Value = '  -0.87%  '
